# Applies row-pair swaps resulting from team name corrections in the
# Indonesia Liga 1 sheet (si 44/45 "Bhayangkara Surabaya United" /
# "Persikabo 1973" swap caused several same-day fixtures to re-sort).
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 18
$ws.Range("B18").Value = 6843497
$ws.Range("F18").Value = 'Persija Jakarta'
$ws.Range("G18").Value = 0
$ws.Range("H18").Value = 0
$ws.Range("I18").Value = 'D'
$ws.Range("J18").Value = 3.25
$ws.Range("K18").Value = 3.25
$ws.Range("L18").Value = 2
$ws.Range("M18").Value = 4
$ws.Range("N18").Value = 3.4
$ws.Range("O18").Value = 1.727
$ws.Range("P18").Value = 0.75
$ws.Range("Q18").Value = 1.75
$ws.Range("R18").Value = 1.95
$ws.Range("S18").Value = 2.5
$ws.Range("T18").Value = 1.9
$ws.Range("U18").Value = 1.9
$ws.Range("W18").Value = 2.4
$ws.Range("X18").Value = -1
$ws.Range("Y18").Value = 0.75
$ws.Range("Z18").Value = -1
$ws.Range("AA18").Value = -1
$ws.Range("AB18").Value = 0.8999999999999999

# Row 19
$ws.Range("B19").Value = 6843498
$ws.Range("F19").Value = 'RANS Nusantara'
$ws.Range("G19").Value = 1
$ws.Range("H19").Value = 2
$ws.Range("I19").Value = 'A'
$ws.Range("J19").Value = 1.285
$ws.Range("K19").Value = 5.25
$ws.Range("L19").Value = 7
$ws.Range("M19").Value = 1.65
$ws.Range("N19").Value = 4.2
$ws.Range("O19").Value = 3.8
$ws.Range("P19").Value = -1
$ws.Range("Q19").Value = 2.075
$ws.Range("R19").Value = 1.725
$ws.Range("S19").Value = 3
$ws.Range("T19").Value = 1.925
$ws.Range("U19").Value = 1.875
$ws.Range("W19").Value = -1
$ws.Range("X19").Value = 2.8
$ws.Range("Y19").Value = -1
$ws.Range("Z19").Value = 0.7250000000000001
$ws.Range("AA19").Value = 0
$ws.Range("AB19").Value = 0

# Row 23
$ws.Range("B23").Value = 6843504
$ws.Range("E23").Value = 'RANS Nusantara'
$ws.Range("F23").Value = 'Persita Tangerang'
$ws.Range("G23").Value = 0
$ws.Range("H23").Value = 1
$ws.Range("I23").Value = 'A'
$ws.Range("J23").Value = 2.1
$ws.Range("L23").Value = 3
$ws.Range("M23").Value = 2.75
$ws.Range("N23").Value = 3.3
$ws.Range("O23").Value = 2.25
$ws.Range("P23").Value = 0.25
$ws.Range("S23").Value = 2.75
$ws.Range("T23").Value = 1.85
$ws.Range("U23").Value = 1.95
$ws.Range("V23").Value = -1
$ws.Range("X23").Value = 1.25
$ws.Range("Y23").Value = -1
$ws.Range("Z23").Value = 1.025
$ws.Range("AA23").Value = -1
$ws.Range("AB23").Value = 0.95

# Row 24
$ws.Range("B24").Value = 6843503
$ws.Range("E24").Value = 'Persik Kediri'
$ws.Range("F24").Value = 'Arema FC'
$ws.Range("G24").Value = 5
$ws.Range("H24").Value = 2
$ws.Range("I24").Value = 'H'
$ws.Range("J24").Value = 1.85
$ws.Range("L24").Value = 3.8
$ws.Range("M24").Value = 2
$ws.Range("N24").Value = 3.1
$ws.Range("O24").Value = 3.4
$ws.Range("P24").Value = -0.25
$ws.Range("S24").Value = 2.5
$ws.Range("T24").Value = 1.975
$ws.Range("U24").Value = 1.825
$ws.Range("V24").Value = 1
$ws.Range("X24").Value = -1
$ws.Range("Y24").Value = 0.7749999999999999
$ws.Range("Z24").Value = -1
$ws.Range("AA24").Value = 0.9750000000000001
$ws.Range("AB24").Value = -1

# Row 25
$ws.Range("B25").Value = 6843505
$ws.Range("E25").Value = 'Persis Solo'
$ws.Range("F25").Value = 'Borneo FC'
$ws.Range("J25").Value = 3
$ws.Range("L25").Value = 2
$ws.Range("M25").Value = 2.8
$ws.Range("O25").Value = 2.1
$ws.Range("P25").Value = 0.25
$ws.Range("Q25").Value = 1.85
$ws.Range("R25").Value = 1.95
$ws.Range("S25").Value = 3
$ws.Range("T25").Value = 1.95
$ws.Range("U25").Value = 1.85
$ws.Range("V25").Value = 1.8
$ws.Range("Y25").Value = 0.8500000000000001
$ws.Range("AA25").Value = 0
$ws.Range("AB25").Value = 0

# Row 26
$ws.Range("B26").Value = 6843506
$ws.Range("E26").Value = 'Bali United'
$ws.Range("F26").Value = 'Madura United'
$ws.Range("J26").Value = 1.833
$ws.Range("L26").Value = 3.5
$ws.Range("M26").Value = 1.85
$ws.Range("O26").Value = 3.5
$ws.Range("P26").Value = -0.5
$ws.Range("Q26").Value = 1.875
$ws.Range("R26").Value = 1.925
$ws.Range("S26").Value = 2.75
$ws.Range("T26").Value = 1.9
$ws.Range("U26").Value = 1.9
$ws.Range("V26").Value = 0.8500000000000001
$ws.Range("Y26").Value = 0.875
$ws.Range("AA26").Value = 0.45
$ws.Range("AB26").Value = -0.5

# Row 52
$ws.Range("B52").Value = 6843533
$ws.Range("E52").Value = 'Borneo FC'
$ws.Range("F52").Value = 'RANS Nusantara'
$ws.Range("H52").Value = 1
$ws.Range("I52").Value = 'D'
$ws.Range("J52").Value = 1.3
$ws.Range("K52").Value = 4.75
$ws.Range("L52").Value = 7.5
$ws.Range("M52").Value = 1.5
$ws.Range("N52").Value = 4.2
$ws.Range("O52").Value = 5
$ws.Range("P52").Value = -1.25
$ws.Range("Q52").Value = 1.95
$ws.Range("R52").Value = 1.75
$ws.Range("S52").Value = 3
$ws.Range("T52").Value = 1.925
$ws.Range("U52").Value = 1.875
$ws.Range("W52").Value = 3.2
$ws.Range("X52").Value = -1
$ws.Range("Z52").Value = 0.75
$ws.Range("AA52").Value = -1
$ws.Range("AB52").Value = 0.875

# Row 53
$ws.Range("B53").Value = 6843532
$ws.Range("E53").Value = 'PSS Sleman'
$ws.Range("F53").Value = 'Persija Jakarta'
$ws.Range("H53").Value = 3
$ws.Range("I53").Value = 'A'
$ws.Range("J53").Value = 3.4
$ws.Range("K53").Value = 3.4
$ws.Range("L53").Value = 1.909
$ws.Range("M53").Value = 2.875
$ws.Range("N53").Value = 3.25
$ws.Range("O53").Value = 2.2
$ws.Range("P53").Value = 0.25
$ws.Range("Q53").Value = 1.775
$ws.Range("R53").Value = 2.025
$ws.Range("S53").Value = 2.5
$ws.Range("T53").Value = 2
$ws.Range("U53").Value = 1.8
$ws.Range("W53").Value = -1
$ws.Range("X53").Value = 1.2
$ws.Range("Z53").Value = 1.025
$ws.Range("AA53").Value = 1
$ws.Range("AB53").Value = -1

# Row 58
$ws.Range("B58").Value = 6843538
$ws.Range("F58").Value = 'Persebaya Surabaya'
$ws.Range("G58").Value = 1
$ws.Range("H58").Value = 2
$ws.Range("M58").Value = 2.2
$ws.Range("N58").Value = 3.3
$ws.Range("Q58").Value = 1.975
$ws.Range("R58").Value = 1.825
$ws.Range("T58").Value = 1.9
$ws.Range("U58").Value = 1.9
$ws.Range("Z58").Value = 0.825
$ws.Range("AA58").Value = 0.45
$ws.Range("AB58").Value = -0.5

# Row 59
$ws.Range("B59").Value = 7030759
$ws.Range("F59").Value = 'PSS Sleman'
$ws.Range("G59").Value = 2
$ws.Range("H59").Value = 3
$ws.Range("M59").Value = 2.15
$ws.Range("N59").Value = 3.5
$ws.Range("Q59").Value = 1.925
$ws.Range("R59").Value = 1.875
$ws.Range("T59").Value = 1.925
$ws.Range("U59").Value = 1.875
$ws.Range("Z59").Value = 0.875
$ws.Range("AA59").Value = 0.925
$ws.Range("AB59").Value = -1

# Row 63
$ws.Range("B63").Value = 6843543
$ws.Range("E63").Value = 'Persija Jakarta'
$ws.Range("F63").Value = 'Borneo FC'
$ws.Range("G63").Value = 1
$ws.Range("I63").Value = 'D'
$ws.Range("J63").Value = 2
$ws.Range("L63").Value = 3.25
$ws.Range("M63").Value = 1.833
$ws.Range("N63").Value = 3.4
$ws.Range("O63").Value = 3.75
$ws.Range("P63").Value = -0.5
$ws.Range("Q63").Value = 1.9
$ws.Range("R63").Value = 1.9
$ws.Range("T63").Value = 1.875
$ws.Range("U63").Value = 1.925
$ws.Range("V63").Value = -1
$ws.Range("W63").Value = 2.4
$ws.Range("Y63").Value = -1
$ws.Range("Z63").Value = 0.8999999999999999
$ws.Range("AA63").Value = -1
$ws.Range("AB63").Value = 0.925

# Row 64
$ws.Range("B64").Value = 6843544
$ws.Range("E64").Value = 'RANS Nusantara'
$ws.Range("F64").Value = 'Madura United'
$ws.Range("G64").Value = 3
$ws.Range("I64").Value = 'H'
$ws.Range("J64").Value = 2.375
$ws.Range("L64").Value = 2.6
$ws.Range("M64").Value = 3
$ws.Range("N64").Value = 3.25
$ws.Range("O64").Value = 2.15
$ws.Range("P64").Value = 0.25
$ws.Range("Q64").Value = 1.85
$ws.Range("R64").Value = 1.95
$ws.Range("T64").Value = 1.825
$ws.Range("U64").Value = 1.975
$ws.Range("V64").Value = 2
$ws.Range("W64").Value = -1
$ws.Range("Y64").Value = 0.8500000000000001
$ws.Range("Z64").Value = -1
$ws.Range("AA64").Value = 0.825
$ws.Range("AB64").Value = -1

# Row 83
$ws.Range("B83").Value = 6843564
$ws.Range("E83").Value = 'Persik Kediri'
$ws.Range("F83").Value = 'PSIS Semarang'
$ws.Range("G83").Value = 1
$ws.Range("H83").Value = 1
$ws.Range("J83").Value = 1.909
$ws.Range("K83").Value = 3.5
$ws.Range("L83").Value = 3.5
$ws.Range("M83").Value = 1.8
$ws.Range("N83").Value = 3.75
$ws.Range("O83").Value = 4
$ws.Range("P83").Value = -0.5
$ws.Range("Q83").Value = 1.8
$ws.Range("R83").Value = 2
$ws.Range("S83").Value = 2.5
$ws.Range("T83").Value = 1.85
$ws.Range("U83").Value = 1.95
$ws.Range("W83").Value = 2.75
$ws.Range("Z83").Value = 1
$ws.Range("AB83").Value = 0.95

# Row 84
$ws.Range("B84").Value = 6843563
$ws.Range("E84").Value = 'Madura United'
$ws.Range("F84").Value = 'Persikabo 1973'
$ws.Range("G84").Value = 0
$ws.Range("H84").Value = 0
$ws.Range("J84").Value = 1.571
$ws.Range("K84").Value = 3.8
$ws.Range("L84").Value = 5
$ws.Range("M84").Value = 1.6
$ws.Range("N84").Value = 4
$ws.Range("O84").Value = 5
$ws.Range("P84").Value = -1
$ws.Range("Q84").Value = 1.95
$ws.Range("R84").Value = 1.85
$ws.Range("S84").Value = 2.75
$ws.Range("T84").Value = 1.8
$ws.Range("U84").Value = 2
$ws.Range("W84").Value = 3
$ws.Range("Z84").Value = 0.8500000000000001
$ws.Range("AB84").Value = 1

# Row 94
$ws.Range("B94").Value = 6843576
$ws.Range("E94").Value = 'RANS Nusantara'
$ws.Range("F94").Value = 'Persik Kediri'
$ws.Range("H94").Value = 0
$ws.Range("I94").Value = 'H'
$ws.Range("J94").Value = 2.1
$ws.Range("L94").Value = 2.9
$ws.Range("M94").Value = 2.1
$ws.Range("N94").Value = 3.3
$ws.Range("O94").Value = 3.1
$ws.Range("P94").Value = -0.25
$ws.Range("T94").Value = 1.95
$ws.Range("U94").Value = 1.85
$ws.Range("V94").Value = 1.1
$ws.Range("W94").Value = -1
$ws.Range("Y94").Value = 0.8500000000000001
$ws.Range("Z94").Value = -1
$ws.Range("AB94").Value = 0.8500000000000001

# Row 95
$ws.Range("B95").Value = 6843575
$ws.Range("E95").Value = 'Persija Jakarta'
$ws.Range("F95").Value = 'Persib Bandung'
$ws.Range("H95").Value = 1
$ws.Range("I95").Value = 'D'
$ws.Range("J95").Value = 1.85
$ws.Range("L95").Value = 3.6
$ws.Range("M95").Value = 1.8
$ws.Range("N95").Value = 3.4
$ws.Range("O95").Value = 3.75
$ws.Range("P95").Value = -0.5
$ws.Range("T95").Value = 1.975
$ws.Range("U95").Value = 1.825
$ws.Range("V95").Value = -1
$ws.Range("W95").Value = 2.4
$ws.Range("Y95").Value = -1
$ws.Range("Z95").Value = 0.95
$ws.Range("AB95").Value = 0.825

# Row 122
$ws.Range("B122").Value = 6843604
$ws.Range("E122").Value = 'Dewa United FC'
$ws.Range("F122").Value = 'Persebaya Surabaya'
$ws.Range("G122").Value = 1
$ws.Range("I122").Value = 'D'
$ws.Range("J122").Value = 2.3
$ws.Range("L122").Value = 2.7
$ws.Range("M122").Value = 2.2
$ws.Range("N122").Value = 3.2
$ws.Range("O122").Value = 2.8
$ws.Range("P122").Value = -0.25
$ws.Range("Q122").Value = 2
$ws.Range("R122").Value = 1.8
$ws.Range("S122").Value = 2.5
$ws.Range("T122").Value = 1.95
$ws.Range("U122").Value = 1.85
$ws.Range("V122").Value = -1
$ws.Range("W122").Value = 2.2
$ws.Range("Y122").Value = -0.5
$ws.Range("Z122").Value = 0.4
$ws.Range("AA122").Value = -1
$ws.Range("AB122").Value = 0.8500000000000001

# Row 123
$ws.Range("B123").Value = 6843602
$ws.Range("E123").Value = 'PSIS Semarang'
$ws.Range("F123").Value = 'PSM Makassar'
$ws.Range("G123").Value = 2
$ws.Range("I123").Value = 'H'
$ws.Range("J123").Value = 2.1
$ws.Range("L123").Value = 3.1
$ws.Range("M123").Value = 1.8
$ws.Range("N123").Value = 3.3
$ws.Range("O123").Value = 3.8
$ws.Range("P123").Value = -0.75
$ws.Range("Q123").Value = 1.875
$ws.Range("R123").Value = 1.925
$ws.Range("S123").Value = 2.25
$ws.Range("T123").Value = 1.85
$ws.Range("U123").Value = 1.95
$ws.Range("V123").Value = 0.8
$ws.Range("W123").Value = -1
$ws.Range("Y123").Value = 0.4375
$ws.Range("Z123").Value = -0.5
$ws.Range("AA123").Value = 0.8500000000000001
$ws.Range("AB123").Value = -1

# Row 188
$ws.Range("B188").Value = 6953283
$ws.Range("E188").Value = 'PSS Sleman'
$ws.Range("F188").Value = 'RANS Nusantara'
$ws.Range("H188").Value = 0
$ws.Range("I188").Value = 'H'
$ws.Range("J188").Value = 2.5
$ws.Range("K188").Value = 3
$ws.Range("L188").Value = 2.6
$ws.Range("M188").Value = 2.2
$ws.Range("N188").Value = 3.1
$ws.Range("O188").Value = 3
$ws.Range("P188").Value = -0.25
$ws.Range("Q188").Value = 1.95
$ws.Range("R188").Value = 1.85
$ws.Range("S188").Value = 2.25
$ws.Range("T188").Value = 1.825
$ws.Range("U188").Value = 1.975
$ws.Range("V188").Value = 1.2
$ws.Range("W188").Value = -1
$ws.Range("Y188").Value = 0.95
$ws.Range("Z188").Value = -1
$ws.Range("AB188").Value = 0.9750000000000001

# Row 189
$ws.Range("B189").Value = 6954176
$ws.Range("E189").Value = 'PSM Makassar'
$ws.Range("F189").Value = 'Persikabo 1973'
$ws.Range("H189").Value = 1
$ws.Range("I189").Value = 'D'
$ws.Range("J189").Value = 1.615
$ws.Range("K189").Value = 3.8
$ws.Range("L189").Value = 4.333
$ws.Range("M189").Value = 1.666
$ws.Range("N189").Value = 4.2
$ws.Range("O189").Value = 3.75
$ws.Range("P189").Value = -0.75
$ws.Range("Q189").Value = 1.825
$ws.Range("R189").Value = 1.975
$ws.Range("S189").Value = 2.75
$ws.Range("T189").Value = 1.85
$ws.Range("U189").Value = 1.95
$ws.Range("V189").Value = -1
$ws.Range("W189").Value = 3.2
$ws.Range("Y189").Value = -1
$ws.Range("Z189").Value = 0.9750000000000001
$ws.Range("AB189").Value = 0.95

# Row 191
$ws.Range("B191").Value = 6953285
$ws.Range("E191").Value = 'Arema FC'
$ws.Range("F191").Value = 'Persis Solo'
$ws.Range("G191").Value = 3
$ws.Range("H191").Value = 1
$ws.Range("J191").Value = 2.6
$ws.Range("K191").Value = 3.2
$ws.Range("L191").Value = 2.4
$ws.Range("M191").Value = 2.7
$ws.Range("N191").Value = 3.2
$ws.Range("O191").Value = 2.3
$ws.Range("P191").Value = 0.25
$ws.Range("Q191").Value = 1.725
$ws.Range("R191").Value = 1.975
$ws.Range("T191").Value = 1.875
$ws.Range("U191").Value = 1.925
$ws.Range("V191").Value = 1.7
$ws.Range("Y191").Value = 0.7250000000000001
$ws.Range("AA191").Value = 0.875
$ws.Range("AB191").Value = -1

# Row 192
$ws.Range("B192").Value = 6953286
$ws.Range("E192").Value = 'Borneo FC'
$ws.Range("F192").Value = 'PSIS Semarang'
$ws.Range("G192").Value = 2
$ws.Range("H192").Value = 0
$ws.Range("J192").Value = 1.6
$ws.Range("K192").Value = 3.6
$ws.Range("L192").Value = 5
$ws.Range("M192").Value = 1.571
$ws.Range("N192").Value = 3.6
$ws.Range("O192").Value = 5.25
$ws.Range("P192").Value = -1
$ws.Range("Q192").Value = 1.9
$ws.Range("R192").Value = 1.9
$ws.Range("T192").Value = 1.825
$ws.Range("U192").Value = 1.975
$ws.Range("V192").Value = 0.571
$ws.Range("Y192").Value = 0.8999999999999999
$ws.Range("AA192").Value = -1
$ws.Range("AB192").Value = 0.9750000000000001

# Row 248
$ws.Range("B248").Value = 6962324
$ws.Range("E248").Value = 'PSM Makassar'
$ws.Range("F248").Value = 'PSS Sleman'
$ws.Range("G248").Value = 2
$ws.Range("H248").Value = 1
$ws.Range("J248").Value = 1.533
$ws.Range("K248").Value = 4
$ws.Range("L248").Value = 4.75
$ws.Range("N248").Value = 4
$ws.Range("O248").Value = 5.5
$ws.Range("P248").Value = -1.25
$ws.Range("Q248").Value = 1.95
$ws.Range("R248").Value = 1.75
$ws.Range("S248").Value = 2.5
$ws.Range("T248").Value = 1.825
$ws.Range("U248").Value = 1.975
$ws.Range("Y248").Value = -0.5
$ws.Range("Z248").Value = 0.375
$ws.Range("AA248").Value = 0.825

# Row 249
$ws.Range("B249").Value = 6962325
$ws.Range("E249").Value = 'Persik Kediri'
$ws.Range("F249").Value = 'RANS Nusantara'
$ws.Range("G249").Value = 4
$ws.Range("H249").Value = 3
$ws.Range("J249").Value = 1.666
$ws.Range("K249").Value = 3.8
$ws.Range("L249").Value = 4
$ws.Range("N249").Value = 4.2
$ws.Range("O249").Value = 5.25
$ws.Range("P249").Value = -1
$ws.Range("Q249").Value = 1.825
$ws.Range("R249").Value = 1.975
$ws.Range("S249").Value = 3
$ws.Range("T249").Value = 1.95
$ws.Range("U249").Value = 1.85
$ws.Range("Y249").Value = 0
$ws.Range("Z249").Value = 0
$ws.Range("AA249").Value = 0.95

# Row 293
$ws.Range("B293").Value = 8056061
$ws.Range("E293").Value = 'Persib Bandung'
$ws.Range("F293").Value = 'Borneo FC'
$ws.Range("G293").Value = 2
$ws.Range("H293").Value = 1
$ws.Range("J293").Value = 1.571
$ws.Range("K293").Value = 3.8
$ws.Range("L293").Value = 4.75
$ws.Range("M293").Value = 1.95
$ws.Range("N293").Value = 3.5
$ws.Range("O293").Value = 3.3
$ws.Range("P293").Value = -0.5
$ws.Range("Q293").Value = 1.975
$ws.Range("R293").Value = 1.825
$ws.Range("S293").Value = 3
$ws.Range("T293").Value = 1.95
$ws.Range("U293").Value = 1.85
$ws.Range("V293").Value = 0.95
$ws.Range("Y293").Value = 0.9750000000000001
$ws.Range("AA293").Value = 0
$ws.Range("AB293").Value = 0

# Row 294
$ws.Range("B294").Value = 8055487
$ws.Range("E294").Value = 'Arema FC'
$ws.Range("F294").Value = 'PSM Makassar'
$ws.Range("G294").Value = 3
$ws.Range("H294").Value = 2
$ws.Range("J294").Value = 1.7
$ws.Range("K294").Value = 3.4
$ws.Range("L294").Value = 4.333
$ws.Range("M294").Value = 2.3
$ws.Range("N294").Value = 3
$ws.Range("O294").Value = 2.9
$ws.Range("P294").Value = -0.25
$ws.Range("Q294").Value = 2.025
$ws.Range("R294").Value = 1.775
$ws.Range("S294").Value = 2.5
$ws.Range("T294").Value = 1.825
$ws.Range("U294").Value = 1.975
$ws.Range("V294").Value = 1.3
$ws.Range("Y294").Value = 1.025
$ws.Range("AA294").Value = 0.825
$ws.Range("AB294").Value = -1

# Row 295
$ws.Range("B295").Value = 8056063
$ws.Range("E295").Value = 'PSIS Semarang'
$ws.Range("F295").Value = 'Bhayangkara Surabaya United'
$ws.Range("G295").Value = 3
$ws.Range("H295").Value = 0
$ws.Range("I295").Value = 'H'
$ws.Range("J295").Value = 1.333
$ws.Range("K295").Value = 4.5
$ws.Range("L295").Value = 7.5
$ws.Range("M295").Value = 1.5
$ws.Range("N295").Value = 4.2
$ws.Range("O295").Value = 5
$ws.Range("P295").Value = -1.25
$ws.Range("Q295").Value = 1.9
$ws.Range("R295").Value = 1.9
$ws.Range("S295").Value = 4
$ws.Range("T295").Value = 1.95
$ws.Range("U295").Value = 1.85
$ws.Range("V295").Value = 0.5
$ws.Range("X295").Value = -1
$ws.Range("Y295").Value = 0.8999999999999999
$ws.Range("Z295").Value = -1
$ws.Range("AA295").Value = -1
$ws.Range("AB295").Value = 0.8500000000000001

# Row 296
$ws.Range("B296").Value = 8056064
$ws.Range("E296").Value = 'Persis Solo'
$ws.Range("F296").Value = 'Persita Tangerang'
$ws.Range("G296").Value = 1
$ws.Range("H296").Value = 2
$ws.Range("I296").Value = 'A'
$ws.Range("J296").Value = 1.85
$ws.Range("K296").Value = 3.5
$ws.Range("L296").Value = 3.4
$ws.Range("M296").Value = 1.85
$ws.Range("N296").Value = 3.8
$ws.Range("O296").Value = 3.4
$ws.Range("P296").Value = -0.5
$ws.Range("Q296").Value = 1.925
$ws.Range("R296").Value = 1.875
$ws.Range("S296").Value = 3.25
$ws.Range("T296").Value = 1.925
$ws.Range("U296").Value = 1.875
$ws.Range("V296").Value = -1
$ws.Range("X296").Value = 2.4
$ws.Range("Y296").Value = -1
$ws.Range("Z296").Value = 0.875
$ws.Range("AA296").Value = -0.5
$ws.Range("AB296").Value = 0.4375

